$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @{ "B"="18.29902407426881"; "D"="3.321983734294835"; "E"="22.22837199877379"; "F"="20.21626514600529"; "G"="3.584377088681463"; "I"="23.46264064641604"; "L"="10.93141291190411"; "O"="17.67226705651923" }
    3 = @{ "B"="17.65096240997742"; "D"="3.303029945144904"; "E"="21.6857031924515"; "F"="20.06091868441453"; "G"="3.587001364298853"; "I"="23.48418673120655"; "L"="10.61887955063239"; "O"="17.6386770513107" }
    4 = @{ "B"="17.24002252170484"; "D"="3.29143081939468"; "E"="21.3477361168436"; "F"="19.97426818049504"; "G"="3.588697700247949"; "I"="23.50656173551024"; "L"="10.42234057075294"; "O"="17.62555077387836" }
    5 = @{ "B"="17.06949089503825"; "D"="3.286716404817432"; "E"="21.20901467006871"; "F"="19.94118974013913"; "G"="3.589410421975019"; "I"="23.51797748065249"; "L"="10.34118407205264"; "O"="17.62208588404232" }
    6 = @{ "B"="17.04099524523132"; "D"="3.285934400318339"; "E"="21.18592594237226"; "F"="19.93583287365052"; "G"="3.589530066584085"; "I"="23.52001174381786"; "L"="10.32764688660641"; "O"="17.62162425064627" }
    7 = @{ "B"="17.2377348164941"; "D"="3.291367185823953"; "E"="21.34586903454753"; "F"="19.97381299047138"; "G"="3.588707225312429"; "I"="23.50670639267062"; "L"="10.42125023873479"; "O"="17.62549641998267" }
    8 = @{ "B"="18.07839291101672"; "D"="3.315442039703176"; "E"="22.04236659460201"; "F"="20.16091269773922"; "G"="3.585264336449044"; "I"="23.46817086866992"; "L"="10.82467380585996"; "O"="17.65912903290348" }
    9 = @{ "B"="19.61560084863116"; "D"="3.362843463261235"; "E"="23.36207995711651"; "F"="20.59524109421572"; "G"="3.579184149170171"; "I"="23.46518209735814"; "L"="11.57485342926473"; "O"="17.78450567751094" }
    10 = @{ "B"="20.66788782884988"; "D"="3.397631026482121"; "E"="24.29327121697127"; "F"="20.95260291744956"; "G"="3.575121676681082"; "I"="23.50713768744078"; "L"="12.09610776561251"; "O"="17.91256504158473" }
    11 = @{ "B"="21.12827698322638"; "D"="3.413415513884966"; "E"="24.7066250366091"; "F"="21.12279103237451"; "G"="3.573360435061992"; "I"="23.53575294890928"; "L"="12.3258508305911"; "O"="17.97850990194764" }
    12 = @{ "B"="21.29987258150918"; "D"="3.41938406147874"; "E"="24.86153348525485"; "F"="21.18827006299184"; "G"="3.57270590526298"; "I"="23.54795116924909"; "L"="12.41172248365569"; "O"="18.00457322947456" }
    13 = @{ "B"="21.26304005779467"; "D"="3.418099060139125"; "E"="24.82824540394294"; "F"="21.17412310164601"; "G"="3.572846318908252"; "I"="23.54526362692032"; "L"="12.39327959325618"; "O"="17.99891173985874" }
    14 = @{ "B"="21.14244983536933"; "D"="3.413906740510841"; "E"="24.71940269947133"; "F"="21.12815771244647"; "G"="3.573306338086357"; "I"="23.53672925548435"; "L"="12.33293848554147"; "O"="17.98063236505383" }
    15 = @{ "B"="21.06822444727526"; "D"="3.411337606534123"; "E"="24.65251835650286"; "F"="21.10013511282332"; "G"="3.573589727821145"; "I"="23.53167883856502"; "L"="12.29582919302693"; "O"="17.96957738695281" }
    16 = @{ "B"="20.63742212245332"; "D"="3.396598412020561"; "E"="24.26603857681531"; "F"="20.94162872938678"; "G"="3.575238518796432"; "I"="23.50545872208314"; "L"="12.08093924270206"; "O"="17.90840896350636" }
    17 = @{ "B"="20.36836688356495"; "D"="3.387543966306685"; "E"="24.02621287129588"; "F"="20.84629808728026"; "G"="3.576272181631497"; "I"="23.49180910223921"; "L"="11.94717228603483"; "O"="17.87284432039864" }
    18 = @{ "B"="20.2118965017548"; "D"="3.382332442620979"; "E"="23.88731462998073"; "F"="20.79218787940311"; "G"="3.57687489093512"; "I"="23.48485576325152"; "L"="11.86954210698524"; "O"="17.85311332892367" }
    19 = @{ "B"="20.15862709604547"; "D"="3.380567371528313"; "E"="23.8401263796626"; "F"="20.77399289783739"; "G"="3.577080363858337"; "I"="23.48265585447735"; "L"="11.84314131498535"; "O"="17.84655764313071" }
    20 = @{ "B"="20.39718682714348"; "D"="3.388508225630372"; "E"="24.0518428515855"; "F"="20.85637198130261"; "G"="3.576161300990329"; "I"="23.4931692858421"; "L"="11.96148401948886"; "O"="17.87655531384286" }
    21 = @{ "B"="21.17794539233754"; "D"="3.415138385151937"; "E"="24.75141749788959"; "F"="21.14163135961537"; "G"="3.573170882919584"; "I"="23.53919911016124"; "L"="12.35069321267586"; "O"="17.98597196656481" }
    22 = @{ "B"="21.67218150690564"; "D"="3.432490733882402"; "E"="25.19913033397475"; "F"="21.3340518058423"; "G"="3.57128880028169"; "I"="23.57721737280778"; "L"="12.59846860888076"; "O"="18.06383535140738" }
    23 = @{ "B"="21.40989884835629"; "D"="3.423235159204453"; "E"="24.96109175791975"; "F"="21.23082722544794"; "G"="3.572286707242905"; "I"="23.55620340387754"; "L"="12.46684978523135"; "O"="18.02170228491368" }
    24 = @{ "B"="20.38416288788476"; "D"="3.388072302442772"; "E"="24.04025870101654"; "F"="20.85181539893507"; "G"="3.576211403823585"; "I"="23.492551561413"; "L"="11.95501594033633"; "O"="17.87487534354424" }
    25 = @{ "B"="19.21269049066636"; "D"="3.350015617516521"; "E"="23.01110935645121"; "F"="20.47080148914001"; "G"="3.580757610522254"; "I"="23.45822308277888"; "L"="11.37683258282155"; "O"="17.74424433659976" }
}

foreach ($row in $data.Keys) {
    $rowData = $data[$row]
    foreach ($col in $rowData.Keys) {
        $ws.Range("$col$row").Value = [double]$rowData[$col]
    }
}

Write-Output "Applied $($data.Count) row updates"